$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '36.589.28'
Set-TextValue $ws.Range('E2') '  +0.55%  '
Set-TextValue $ws.Range('D3') '2.005.22'
Set-TextValue $ws.Range('E3') '  -0.45%  '
Set-TextValue $ws.Range('E4') '  +0.02%  '
Set-TextValue $ws.Range('D5') '247.18'
Set-TextValue $ws.Range('E5') '  -2.06%  '
Set-TextValue $ws.Range('E6') '  -1.40%  '
Set-TextValue $ws.Range('D7') '62.64'
Set-TextValue $ws.Range('E7') '  +0.66%  '
Set-TextValue $ws.Range('E8') '  +0.06%  '
Set-TextValue $ws.Range('D9') '0.385'
Set-TextValue $ws.Range('E9') '  +3.86%  '
Set-TextValue $ws.Range('D10') '57.21'
Set-TextValue $ws.Range('E10') '  -1.93%  '
Set-TextValue $ws.Range('D11') '0.0780'
Set-TextValue $ws.Range('E11') '  +4.95%  '
Set-TextValue $ws.Range('D12') '0.104'
Set-TextValue $ws.Range('E12') '  -0.16%  '
Set-TextValue $ws.Range('D13') '0.889'
Set-TextValue $ws.Range('E13') '  -1.92%  '
Set-TextValue $ws.Range('D14') '22.56'
Set-TextValue $ws.Range('E14') '  +9.06%  '
Set-TextValue $ws.Range('D15') '14.15'
Set-TextValue $ws.Range('E15') '  -5.34%  '
Set-TextValue $ws.Range('D16') '2.302.94'
Set-TextValue $ws.Range('E16') '  -0.21%  '
Set-TextValue $ws.Range('D17') '5.50'
Set-TextValue $ws.Range('E17') '  +0.11%  '
Set-TextValue $ws.Range('D18') '2.008.72'
Set-TextValue $ws.Range('E18') '  -0.34%  '
Set-TextValue $ws.Range('D19') '36.521.82'
Set-TextValue $ws.Range('E19') '  +0.47%  '
Set-TextValue $ws.Range('E20') '  -0.13%  '
Set-TextValue $ws.Range('D21') '0.0₃0871'
Set-TextValue $ws.Range('E21') '  +0.58%  '
Set-TextValue $ws.Range('D22') '5.31'
Set-TextValue $ws.Range('E22') '  +0.16%  '
Set-TextValue $ws.Range('D23') '238.26'
Set-TextValue $ws.Range('E23') '  +1.61%  '
Set-TextValue $ws.Range('E24') '  -0.06%  '
Set-TextValue $ws.Range('D25') '2.52'
Set-TextValue $ws.Range('E25') '  -7.31%  '
Set-TextValue $ws.Range('E26') '  +0.41%  '
Set-TextValue $ws.Range('D27') '9.90'
Set-TextValue $ws.Range('E27') '  +2.69%  '
Set-TextValue $ws.Range('E28') '  +26.39%  '
Set-TextValue $ws.Range('D29') '160.06'
Set-TextValue $ws.Range('E29') '  -2.01%  '
Set-TextValue $ws.Range('D30') '20.14'
Set-TextValue $ws.Range('E30') '  +2.50%  '
Set-TextValue $ws.Range('E31') '  +0.73%  '
Set-TextValue $ws.Range('E32') '  -0.43%  '
Set-TextValue $ws.Range('D33') '5.01'
Set-TextValue $ws.Range('E33') '  -2.47%  '
Set-TextValue $ws.Range('D34') '0.0624'
Set-TextValue $ws.Range('E34') '  +2.52%  '
Set-TextValue $ws.Range('E35') '  -2.22%  '
Set-TextValue $ws.Range('D36') '6.52'
Set-TextValue $ws.Range('E36') '  +10.06%  '
Set-TextValue $ws.Range('E37') '  -3.57%  '
Set-TextValue $ws.Range('E38') '  +0.06%  '
Set-TextValue $ws.Range('E39') '  +0.81%  '
Set-TextValue $ws.Range('E40') '  +20.03%  '
Set-TextValue $ws.Range('D41') '1.27'
Set-TextValue $ws.Range('E41') '  +3.26%  '
Set-TextValue $ws.Range('E42') '  -4.11%  '
Set-TextValue $ws.Range('D43') '2.92'
Set-TextValue $ws.Range('E43') '  +0.04%  '
Set-TextValue $ws.Range('E44') '  -0.79%  '
Set-TextValue $ws.Range('E45') '  -1.18%  '
Set-TextValue $ws.Range('D46') '16.76'
Set-TextValue $ws.Range('E46') '  -2.17%  '
Set-TextValue $ws.Range('D47') '93.37'
Set-TextValue $ws.Range('E47') '  -1.86%  '
Set-TextValue $ws.Range('D48') '7.66'
Set-TextValue $ws.Range('E48') '  -4.94%  '
Set-TextValue $ws.Range('D49') '1.359.22'
Set-TextValue $ws.Range('E49') '  -6.71%  '
Set-TextValue $ws.Range('D50') '2.88'
Set-TextValue $ws.Range('E50') '  -1.74%  '
Set-TextValue $ws.Range('D51') '2.195.13'
